$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.642.12"
$ws.Range("E2").Value = "  -0.24%  "

$ws.Range("D3").Value = "1.642.65"
$ws.Range("E3").Value = "  +0.54%  "

$ws.Range("D5").Value = "215.38"
$ws.Range("E5").Value = "  +0.79%  "

$ws.Range("E6").Value = "  +1.12%  "

$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("D9").Value = "0.0628"
$ws.Range("E9").Value = "  +0.77%  "

$ws.Range("D10").Value = "19.30"
$ws.Range("E10").Value = "  +0.37%  "

$ws.Range("D11").Value = "0.0841"
$ws.Range("E11").Value = "  -0.07%  "

$ws.Range("D12").Value = "1.870.90"

$ws.Range("D13").Value = "1.647.30"
$ws.Range("E13").Value = "  +0.48%  "

$ws.Range("E14").Value = "  +2.25%  "

$ws.Range("E15").Value = "  +1.10%  "

$ws.Range("D16").Value = "65.35"
$ws.Range("E16").Value = "  +2.75%  "

$ws.Range("D17").Value = "26.688.08"
$ws.Range("E17").Value = "  -0.02%  "

$ws.Range("D18").Value = "0.0₃0745"
$ws.Range("E18").Value = "  +0.36%  "

$ws.Range("D19").Value = "217.09"
$ws.Range("E19").Value = "  -0.79%  "

$ws.Range("E21").Value = "  +1.11%  "

$ws.Range("E22").Value = "  +2.54%  "

$ws.Range("D23").Value = "9.51"
$ws.Range("E23").Value = "  +1.51%  "

$ws.Range("D24").Value = "2.23"
$ws.Range("E24").Value = "  +13.70%  "

$ws.Range("D25").Value = "145.54"
$ws.Range("E25").Value = "  -1.70%  "

$ws.Range("E26").Value = "  +0.18%  "

$ws.Range("E27").Value = "  -1.00%  "

$ws.Range("D28").Value = "7.17"
$ws.Range("E28").Value = "  +4.28%  "

$ws.Range("E29").Value = "  +1.39%  "

$ws.Range("E30").Value = "  +2.56%  "

$ws.Range("E31").Value = "  +0.76%  "

$ws.Range("E32").Value = "  +2.35%  "

$ws.Range("D33").Value = "3.05"
$ws.Range("E33").Value = "  +1.97%  "

$ws.Range("D34").Value = "1.277.42"
$ws.Range("E34").Value = "  +4.20%  "

$ws.Range("E35").Value = "  +2.64%  "

$ws.Range("E36").Value = "  +4.78%  "

$ws.Range("E37").Value = "  +0.35%  "

$ws.Range("D38").Value = "0.535"
$ws.Range("E38").Value = "  +6.61%  "

$ws.Range("E39").Value = "  +2.58%  "

$ws.Range("E40").Value = "  +0.12%  "

$ws.Range("D41").Value = "0.818"
$ws.Range("E41").Value = "  +2.77%  "

$ws.Range("E42").Value = "  -1.49%  "

$ws.Range("E43").Value = "  +1.91%  "

$ws.Range("D44").Value = "1.781.50"
$ws.Range("E44").Value = "  +0.69%  "

$ws.Range("D45").Value = "92.17"
$ws.Range("E45").Value = "  -0.56%  "

$ws.Range("D46").Value = "59.92"
$ws.Range("E46").Value = "  +8.14%  "

$ws.Range("D47").Value = "1.59"
$ws.Range("E47").Value = "  +1.76%  "

$ws.Range("E48").Value = "  +0.64%  "

$ws.Range("D49").Value = "7.80"
$ws.Range("E49").Value = "  +1.67%  "

$ws.Range("D50").Value = "0.0971"
$ws.Range("E50").Value = "  +2.95%  "

$ws.Range("E51").Value = "  -0.74%  "
